$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 168.16667
$ws.Range("I2").Value = 135
$ws.Range("J2").Value = 201.33333
$ws.Range("K2").Value = 135
$ws.Range("L2").Value = 201.33333
$ws.Range("M2").Value = -22
$ws.Range("N2").Value = -427.33333

$ws.Range("H19").Value = 520.7037
$ws.Range("I19").Value = 667.2143
$ws.Range("J19").Value = 362.92307
$ws.Range("K19").Value = 667.2143
$ws.Range("L19").Value = 362.92307
$ws.Range("M19").Value = -492.2143
$ws.Range("N19").Value = -712.9230700000001

$ws.Range("H21").Value = 50301.43
$ws.Range("I21").Value = 55351.668
$ws.Range("J21").Value = 20000
$ws.Range("K21").Value = 55351.668
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = -54883.668
$ws.Range("N21").Value = -20936

$ws.Range("H23").Value = 50301.43
$ws.Range("I23").Value = 55351.668
$ws.Range("J23").Value = 20000
$ws.Range("K23").Value = 55351.668
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = -55117.668
$ws.Range("N23").Value = -20468

$ws.Range("H26").Value = 800
$ws.Range("I26").Value = 800
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 800
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -456

$ws.Range("H29").Value = 279.18182
$ws.Range("I29").Value = 279.18182
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 837.54546
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -556.54546

$ws.Range("H32").Value = 841
$ws.Range("I32").Value = 850
$ws.Range("J32").Value = 836.5
$ws.Range("K32").Value = 850
$ws.Range("L32").Value = 836.5
$ws.Range("M32").Value = -524
$ws.Range("N32").Value = -1488.5

$ws.Range("H38").Value = 620.2174
$ws.Range("I38").Value = 48.333332
$ws.Range("J38").Value = 987.8570999999999
$ws.Range("K38").Value = 144.999996
$ws.Range("L38").Value = 2963.5713
$ws.Range("M38").Value = 227.000004
$ws.Range("N38").Value = -3707.5713

$ws.Range("H42").Value = 65.09999999999999
$ws.Range("I42").Value = 33.25
$ws.Range("J42").Value = 86.333336
$ws.Range("K42").Value = 99.75
$ws.Range("L42").Value = 259.000008
$ws.Range("M42").Value = 130.25
$ws.Range("N42").Value = -719.000008

$ws.Range("H43").Value = 1540.8846
$ws.Range("I43").Value = 1658.3334
$ws.Range("J43").Value = 1440.2142
$ws.Range("K43").Value = 1658.3334
$ws.Range("L43").Value = 1440.2142
$ws.Range("M43").Value = -1589.3334
$ws.Range("N43").Value = -1578.2142

$ws.Range("H53").Value = 274.6
$ws.Range("I53").Value = 189.33333
$ws.Range("J53").Value = 402.5
$ws.Range("K53").Value = 189.33333
$ws.Range("L53").Value = 402.5
$ws.Range("M53").Value = 447.66667
$ws.Range("N53").Value = -1676.5

$ws.Range("H58").Value = 597.5909
$ws.Range("I58").Value = 211.18182
$ws.Range("J58").Value = 984
$ws.Range("K58").Value = 633.5454599999999
$ws.Range("L58").Value = 2952
$ws.Range("M58").Value = -483.5454599999999
$ws.Range("N58").Value = -3252

$ws.Range("H87").Value = 40000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 40000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42496

$ws.Range("H90").Value = 40000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 40000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -132480

$ws.Range("H128").Value = 41226.668
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 41226.668
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 41226.668
$ws.Range("N128").Value = -51186.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 708
$ws.Range("I39").Value = 708
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 708
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -188

$ws.Range("H41").Value = 1345.3334
$ws.Range("I41").Value = 1345.3334
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1345.3334
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -931.3334

$ws.Range("H59").Value = 24000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 24000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 24000
$ws.Range("N59").Value = -25608

$ws.Range("H102").Value = 184349.83
$ws.Range("I102").Value = 251013
$ws.Range("J102").Value = 51023.5
$ws.Range("K102").Value = 251013
$ws.Range("L102").Value = 51023.5
$ws.Range("M102").Value = -249391
$ws.Range("N102").Value = -54267.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1698.3334
$ws.Range("I99").Value = 1166.6666
$ws.Range("J99").Value = 1875.5555
$ws.Range("K99").Value = 1166.6666
$ws.Range("L99").Value = 1875.5555
$ws.Range("M99").Value = 331.3334
$ws.Range("N99").Value = -4871.5555

$ws.Range("H105").Value = 2188.5
$ws.Range("I105").Value = 1561.8334
$ws.Range("J105").Value = 2815.1667
$ws.Range("K105").Value = 1561.8334
$ws.Range("L105").Value = 2815.1667
$ws.Range("M105").Value = 185.1666
$ws.Range("N105").Value = -6309.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 10133.333
$ws.Range("I45").Value = 6000
$ws.Range("J45").Value = 10960
$ws.Range("K45").Value = 6000
$ws.Range("L45").Value = 10960
$ws.Range("M45").Value = -5407
$ws.Range("N45").Value = -12146

$ws.Range("H54").Value = 9000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 9000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 9000
$ws.Range("N54").Value = -10316

$ws.Range("H99").Value = 1548.5555
$ws.Range("I99").Value = 1366.1538
$ws.Range("J99").Value = 2022.8
$ws.Range("K99").Value = 1366.1538
$ws.Range("L99").Value = 2022.8
$ws.Range("M99").Value = 131.8462
$ws.Range("N99").Value = -5018.8

$ws.Range("H105").Value = 5376.375
$ws.Range("I105").Value = 4282.2
$ws.Range("J105").Value = 7200
$ws.Range("K105").Value = 4282.2
$ws.Range("L105").Value = 7200
$ws.Range("M105").Value = -2535.2
$ws.Range("N105").Value = -10694

$ws.Range("H122").Value = 1110.069
$ws.Range("I122").Value = 1018.5333
$ws.Range("J122").Value = 1208.1428
$ws.Range("K122").Value = 3055.5999
$ws.Range("L122").Value = 3624.4284
$ws.Range("M122").Value = -605.5999000000002
$ws.Range("N122").Value = -8524.428400000001

$ws.Range("H126").Value = 1548.5555
$ws.Range("I126").Value = 1366.1538
$ws.Range("J126").Value = 2022.8
$ws.Range("K126").Value = 4098.4614
$ws.Range("L126").Value = 6068.4
$ws.Range("M126").Value = -1628.4614
$ws.Range("N126").Value = -11008.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 2015
$ws.Range("I35").Value = 2015
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2015
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1717

$ws.Range("H47").Value = 4500
$ws.Range("I47").Value = 4000
$ws.Range("J47").Value = 5000
$ws.Range("K47").Value = 4000
$ws.Range("L47").Value = 5000
$ws.Range("M47").Value = -3432
$ws.Range("N47").Value = -6136

$ws.Range("H48").Value = 14666.667
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 14666.667
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 14666.667
$ws.Range("N48").Value = -15636.667

$ws.Range("H52").Value = 7515
$ws.Range("I52").Value = 3030
$ws.Range("J52").Value = 12000
$ws.Range("K52").Value = 3030
$ws.Range("L52").Value = 12000
$ws.Range("M52").Value = -2771
$ws.Range("N52").Value = -12518

$ws.Range("H55").Value = 3750
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 9000
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 9000
$ws.Range("M55").Value = -1673
$ws.Range("N55").Value = -9654

$ws.Range("H102").Value = 2985.7097
$ws.Range("I102").Value = 3085.0454
$ws.Range("J102").Value = 2742.889
$ws.Range("K102").Value = 3085.0454
$ws.Range("L102").Value = 2742.889
$ws.Range("M102").Value = -1463.0454
$ws.Range("N102").Value = -5986.889

$ws.Range("H134").Value = 29775.334
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 29775.334
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 89326.00199999999
$ws.Range("N134").Value = -94396.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8286.462
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 9702.182000000001
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 9702.182000000001
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -10292.182

$ws.Range("H27").Value = 8286.462
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 9702.182000000001
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 9702.182000000001
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -9916.182000000001

$ws.Range("H61").Value = 1779.8
$ws.Range("I61").Value = 1399
$ws.Range("J61").Value = 2215
$ws.Range("K61").Value = 1399
$ws.Range("L61").Value = 2215
$ws.Range("M61").Value = -1197
$ws.Range("N61").Value = -2619

$ws.Range("H100").Value = 2338.2666
$ws.Range("I100").Value = 2595.7144
$ws.Range("J100").Value = 2113
$ws.Range("K100").Value = 2595.7144
$ws.Range("L100").Value = 2113
$ws.Range("M100").Value = -2054.7144
$ws.Range("N100").Value = -3195

$ws.Range("H113").Value = 1779.8
$ws.Range("I113").Value = 1399
$ws.Range("J113").Value = 2215
$ws.Range("K113").Value = 1399
$ws.Range("L113").Value = 2215
$ws.Range("M113").Value = 771
$ws.Range("N113").Value = -6555

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 4250
$ws.Range("I38").Value = 3000
$ws.Range("J38").Value = 4666.6665
$ws.Range("K38").Value = 3000
$ws.Range("L38").Value = 4666.6665
$ws.Range("M38").Value = -2527
$ws.Range("N38").Value = -5612.6665

$ws.Range("H47").Value = 6000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 6000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 6000
$ws.Range("N47").Value = -7144

$ws.Range("H48").Value = 8000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 8000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 8000
$ws.Range("N48").Value = -9138

$ws.Range("H49").Value = 7000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 7000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 7000
$ws.Range("N49").Value = -7460
$ws.Range("M49").ClearContents()

$ws.Range("H132").Value = 3603.0852
$ws.Range("I132").Value = 5016.5356
$ws.Range("J132").Value = 1520.1052
$ws.Range("K132").Value = 15049.6068
$ws.Range("L132").Value = 4560.3156
$ws.Range("M132").Value = -12519.6068
$ws.Range("N132").Value = -9620.3156
